$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").Value = "61.371.52"
$ws.Range("E2").Value = "  -6.18%  "

# Row 3
$ws.Range("D3").Value = "3.117.00"
$ws.Range("E3").Value = "  -7.80%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Formula = "'507.58"
$ws.Range("E5").Value = "  -4.02%  "

# Row 6
$ws.Range("D6").Formula = "'165.89"
$ws.Range("E6").Value = "  -11.68%  "

# Row 7
$ws.Range("D7").Formula = "'0.578"
$ws.Range("E7").Value = "  -4.40%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").Value = "3.117.15"
$ws.Range("E9").Value = "  -7.64%  "

# Row 10
$ws.Range("D10").Formula = "'0.580"
$ws.Range("E10").Value = "  -7.52%  "

# Row 11
$ws.Range("D11").Formula = "'51.22"
$ws.Range("E11").Value = "  -12.66%  "

# Row 12
$ws.Range("D12").Formula = "'0.125"
$ws.Range("E12").Value = "  -6.38%  "

# Row 13
$ws.Range("D13").Formula = "'0.0000242"
$ws.Range("E13").Value = "  -5.04%  "

# Row 14
$ws.Range("D14").Formula = "'8.61"
$ws.Range("E14").Value = "  -7.04%  "

# Row 15
$ws.Range("D15").Value = "3.625.27"
$ws.Range("E15").Value = "  -7.81%  "

# Row 16
$ws.Range("E16").Value = "  -8.92%  "

# Row 17
$ws.Range("D17").Value = "3.125.12"
$ws.Range("E17").Value = "  -7.87%  "

# Row 18
$ws.Range("D18").Value = "61.351.77"
$ws.Range("E18").Value = "  -5.96%  "

# Row 19
$ws.Range("D19").Formula = "'16.61"
$ws.Range("E19").Value = "  -4.92%  "

# Row 20
$ws.Range("D20").Formula = "'10.48"
$ws.Range("E20").Value = "  -5.78%  "

# Row 21
$ws.Range("D21").Formula = "'0.929"
$ws.Range("E21").Value = "  -4.50%  "

# Row 22
$ws.Range("D22").Formula = "'352.40"
$ws.Range("E22").Value = "  -5.33%  "

# Row 23
$ws.Range("D23").Formula = "'78.46"
$ws.Range("E23").Value = "  -3.68%  "

# Row 24
$ws.Range("D24").Formula = "'3.60"
$ws.Range("E24").Value = "  -3.75%  "

# Row 25
$ws.Range("D25").Formula = "'10.58"
$ws.Range("E25").Value = "  -2.45%  "

# Row 26
$ws.Range("D26").Formula = "'6.07"
$ws.Range("E26").Value = "  +4.24%  "

# Row 27
$ws.Range("E27").Value = "  +1.42%  "

# Row 28
$ws.Range("D28").Formula = "'2.51"
$ws.Range("E28").Value = "  -5.23%  "

# Row 29
$ws.Range("D29").Formula = "'10.73"
$ws.Range("E29").Value = "  -7.12%  "

# Row 30
$ws.Range("D30").Formula = "'7.84"
$ws.Range("E30").Value = "  -8.14%  "

# Row 31
$ws.Range("D31").Formula = "'631.23"
$ws.Range("E31").Value = "  -6.30%  "

# Row 32
$ws.Range("D32").Formula = "'27.32"
$ws.Range("E32").Value = "  -7.88%  "

# Row 33
$ws.Range("D33").Formula = "'6.15"
$ws.Range("E33").Value = "  -8.25%  "

# Row 34
$ws.Range("D34").Formula = "'10.87"
$ws.Range("E34").Value = "  -2.80%  "

# Row 35
$ws.Range("E35").Value = "  -0.03%  "

# Row 36
$ws.Range("D36").Formula = "'0.0997"
$ws.Range("E36").Value = "  -5.77%  "

# Row 37
$ws.Range("D37").Formula = "'55.11"
$ws.Range("E37").Value = "  -9.83%  "

# Row 38
$ws.Range("D38").Formula = "'34.89"
$ws.Range("E38").Value = "  -4.63%  "

# Row 39
$ws.Range("D39").Formula = "'0.360"
$ws.Range("E39").Value = "  -5.60%  "

# Row 40
$ws.Range("D40").Formula = "'1.00"
$ws.Range("E40").Value = "  +0.17%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0658"
$ws.Range("E41").Value = "  +5.93%  "

# Row 42
$ws.Range("D42").Formula = "'0.118"
$ws.Range("E42").Value = "  -7.05%  "

# Row 43
$ws.Range("D43").Value = "2.777.75"
$ws.Range("E43").Value = "  -2.88%  "

# Row 44
$ws.Range("D44").Formula = "'2.40"
$ws.Range("E44").Value = "  +1.84%  "

# Row 45
$ws.Range("D45").Formula = "'2.58"
$ws.Range("E45").Value = "  -2.92%  "

# Row 46
$ws.Range("D46").Formula = "'2.78"
$ws.Range("E46").Value = "  +6.86%  "

# Row 47
$ws.Range("D47").Formula = "'0.0372"
$ws.Range("E47").Value = "  -5.95%  "

# Row 48
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Formula = "'2.45"
$ws.Range("E48").Value = "  -11.32%  "

# Row 49
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Formula = "'2.88"
$ws.Range("E49").Value = "  +1.28%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Formula = "'0.119"
$ws.Range("E50").Value = "  -5.03%  "

# Row 51
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Formula = "'129.90"
$ws.Range("E51").Value = "  -6.00%  "
